$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two added columns (C, D)
$ws.Range("C1").Value = "128 Channels, 45° FOV"
$ws.Range("D1").Value = "128 Channels, 22.5° FOV"

# Match the header formatting used by the existing "64 Channels, 45° FOV" header (B1)
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)  # xlPasteFormats

# New data columns: C = "128 Channels, 45° FOV", D = "128 Channels, 22.5° FOV"
$cValues = @(0, 6, 4, 3, 1, 1, 1, 0)
$dValues = @(0, 13, 9, 5, 3, 1, 2, 0)

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
}
